# Outstandings.xlsx - "Add files via upload"
#
# Adds a new outstanding-invoice line (I-C-1-23-453322, 25-Oct-2023,
# 52510) for "Collective Trade Links Pvt Ltd" on the "Purchase 22-23"
# sheet, right under the two existing invoices for that vendor, and makes
# that sheet the active/selected one (selection D28) instead of "Sale 22-23".

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Purchase 22-23"
$ws2 = $wb.Worksheets.Item(2)   # "Sale 22-23"

# Insert a new row 6 - this shifts the old rows 7/9 down to 8/10 and
# renumbers them automatically, matching the diff.
$ws1.Rows.Item(6).Insert()

# Copy the formatting of the row above (row 5, same vendor block) onto the
# freshly inserted row so fonts/borders/number formats line up.
$ws1.Range("A5:F5").Copy()
$ws1.Range("A6:F6").PasteSpecial(-4122)   # -4122 = xlPasteFormats
$ws1.Application.CutCopyMode = $false

# Fill in the new invoice's data.
$ws1.Range("A6").Value = ""
$ws1.Range("B6").Value = 45224
$ws1.Range("C6").Value = "I-C-1-23-453322"
$ws1.Range("D6").Value = "Collective Trade Links Pvt Ltd"
$ws1.Range("E6").Value = 52510
$ws1.Range("F6").Formula = "=E4+E5+E6"

# The running-total formula moves from the old F5 down to the new F6, so
# F5 goes back to being a plain (empty) cell.
$ws1.Range("F5").Formula = ""

# Make "Purchase 22-23" the active sheet with D28 selected (was "Sale
# 22-23" with tabSelected before).
$ws1.Activate()
$ws1.Range("D28").Select()
